$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $rest = $parts[1..($parts.Count - 1)]
            $newParts = $rest + @("System")
            $newVal = [string]::Join(", ", $newParts)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
